$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (unchanged text, just ensure values) ---
$ws.Range("A1").Value = "Polymer"
$ws.Range("B1").Value = "Effective Diameter Avg. (nm)"
$ws.Range("C1").Value = "Std Diam. "
$ws.Range("D1").Value = "PDI"
$ws.Range("E1").Value = "Std PDI"

# --- Data rows, reordered (row 6 label corrected DMA B2 -> DMA B1) ---
$rows = @(
    @{ Row = 2;  Label = "DIP S1"; Start = 17; End = 19 },
    @{ Row = 3;  Label = "DMA S1"; Start = 20; End = 22 },
    @{ Row = 4;  Label = "DIP B1"; Start = 11; End = 13 },
    @{ Row = 5;  Label = "DMA B1"; Start = 14; End = 16 },
    @{ Row = 6;  Label = "DIP G2"; Start = 8;  End = 10 },
    @{ Row = 7;  Label = "DMA G1"; Start = 5;  End = 7  },
    @{ Row = 8;  Label = "DMA G2"; Start = 2;  End = 4  }
)

foreach ($r in $rows) {
    $row = $r.Row
    $start = $r.Start
    $end = $r.End

    $ws.Cells.Item($row, 1).Value = $r.Label
    $ws.Cells.Item($row, 2).Formula = "=AVERAGE('[1]20231016_DLSData'!`$G`$$start`:`$G`$$end)"
    $ws.Cells.Item($row, 3).Formula = "=STDEV('[2]20231016_DLSData'!`$G`$$start`:`$G`$$end)"
    $ws.Cells.Item($row, 4).Formula = "=AVERAGE('[2]20231016_DLSData'!`$H`$$start`:`$H`$$end)"
    $ws.Cells.Item($row, 5).Formula = "=STDEV('[2]20231016_DLSData'!`$H`$$start`:`$H`$$end)"
}

# --- Selection change ---
$ws.Range("B1:E1").Select()

# --- Window view change ---
$excel.ActiveWindow.Left = 3420
$excel.ActiveWindow.Top = 3400
$excel.ActiveWindow.Width = 18600
$excel.ActiveWindow.Height = 13240
